# Fix the YouTube / Soundcloud / 1001Tracklists collection data.
# The old code used a "merge" instead of a "join", which silently dropped
# the new song ("Martin Garrix, Tove Lo" - "Pressure" on "STMPD RCRDS")
# from the dataset and left every other metric under-counted. This script
# rewrites each summary sheet with the corrected, complete numbers.

$wb = $excel.ActiveWorkbook

function Set-SheetRows {
    # NOTE: positional parameters only - named parameter binding with
    # array arguments is unreliable in this host, so always call this
    # function positionally: Set-SheetRows $ws $labelStyleCol $rows
    param($ws, $labelStyleCol, $rows)

    $r = 2
    foreach ($row in $rows) {
        $c = 1
        foreach ($val in $row) {
            $cell = $ws.Cells.Item($r, $c)
            $cell.Value = $val
            if ($labelStyleCol -gt 0 -and $c -eq $labelStyleCol) {
                # Re-apply the bordered/bold/centered header style (style
                # index 1 in styles.xml) explicitly: assigning .Value above
                # resets a cell back to the default style, and re-assigning
                # a previously captured .Style object does not stick in
                # this host, so set the individual format properties instead
                # (this reproduces the existing style so Excel reuses it
                # rather than creating a duplicate).
                $cell.Font.Bold = $true
                $cell.HorizontalAlignment = -4108
                $cell.VerticalAlignment = -4160
                $cell.Borders.LineStyle = 1
            }
            $c = $c + 1
        }
        $r = $r + 1
    }
}

# ---------------------------------------------------------------------
# By_Track_YouTube
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("By_Track_YouTube")
$rows1 = @(
    ,@("Martin Garrix, Tove Lo", "Pressure", "STMPD RCRDS", 1286202)
    ,@("Feint", "Do Better", "Monstercat", 63216)
    ,@("The Bloody Beetroots, Teddy Killerz", "Elevate", "Monstercat", 57445)
    ,@("Bleu Clair, OOTORO", "Beat Like This", "STMPD RCRDS", 40680)
)
Set-SheetRows $ws1 0 $rows1

# ---------------------------------------------------------------------
# By_Track_1001Tracklists
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("By_Track_1001Tracklists")
$rows2 = @(
    ,@("Bleu Clair, OOTORO", "Beat Like This", "STMPD RCRDS", 41, 44)
    ,@("Martin Garrix, Tove Lo", "Pressure", "STMPD RCRDS", 18, 19)
    ,@("The Bloody Beetroots, Teddy Killerz", "Elevate", "Monstercat", 1, 3)
    ,@("Feint", "Do Better", "Monstercat", 0, 0)
)
Set-SheetRows $ws2 0 $rows2

# ---------------------------------------------------------------------
# By_Track_Soundcloud
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("By_Track_Soundcloud")
$rows3 = @(
    ,@("Martin Garrix, Tove Lo", "Pressure", "STMPD RCRDS", 59108)
    ,@("Feint", "Do Better", "Monstercat", 33805)
    ,@("Bleu Clair, OOTORO", "Beat Like This", "STMPD RCRDS", 32523)
    ,@("The Bloody Beetroots, Teddy Killerz", "Elevate", "Monstercat", 25372)
)
Set-SheetRows $ws3 0 $rows3

# ---------------------------------------------------------------------
# By_Artist_YouTube
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("By_Artist_YouTube")
$rows4 = @(
    ,@("Martin Garrix", 1286202)
    ,@("Tove Lo", 1286202)
    ,@("Feint", 63216)
    ,@("Teddy Killerz", 57445)
    ,@("The Bloody Beetroots", 57445)
    ,@("Bleu Clair", 40680)
    ,@("OOTORO", 40680)
)
Set-SheetRows $ws4 1 $rows4

# ---------------------------------------------------------------------
# By_Artist_1001Tracklists
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("By_Artist_1001Tracklists")
$rows5 = @(
    ,@("Bleu Clair", 41, 44)
    ,@("OOTORO", 41, 44)
    ,@("Martin Garrix", 18, 19)
    ,@("Tove Lo", 18, 19)
    ,@("Teddy Killerz", 1, 3)
    ,@("The Bloody Beetroots", 1, 3)
    ,@("Feint", 0, 0)
)
Set-SheetRows $ws5 1 $rows5

# ---------------------------------------------------------------------
# By_Artist_Soundcloud
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("By_Artist_Soundcloud")
$rows6 = @(
    ,@("Martin Garrix", 59108)
    ,@("Tove Lo", 59108)
    ,@("Feint", 33805)
    ,@("Bleu Clair", 32523)
    ,@("OOTORO", 32523)
    ,@("Teddy Killerz", 25372)
    ,@("The Bloody Beetroots", 25372)
)
Set-SheetRows $ws6 1 $rows6

# ---------------------------------------------------------------------
# By_Label_YouTube
# ---------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item("By_Label_YouTube")
$rows7 = @(
    ,@("STMPD RCRDS", 1326882)
    ,@("Monstercat", 120661)
)
Set-SheetRows $ws7 0 $rows7

# ---------------------------------------------------------------------
# By_Label_1001Tracklists
# ---------------------------------------------------------------------
$ws8 = $wb.Worksheets.Item("By_Label_1001Tracklists")
$rows8 = @(
    ,@("STMPD RCRDS", 59, 63)
    ,@("Monstercat", 1, 3)
)
Set-SheetRows $ws8 0 $rows8

# ---------------------------------------------------------------------
# By_Label_Soundcloud
# ---------------------------------------------------------------------
$ws9 = $wb.Worksheets.Item("By_Label_Soundcloud")
$rows9 = @(
    ,@("STMPD RCRDS", 91631)
    ,@("Monstercat", 59177)
)
Set-SheetRows $ws9 0 $rows9

Write-Host "Updated all 9 summary sheets with corrected join-based totals."
